$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# Fill in the new test case row (row 28)
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "testcase27_csv_csv_bigdata_match"
$ws.Range("C28").Formula = '=CONCAT("/app/test/testcases/",B28,".xlsx")'
$ws.Range("D28").Value = "N"

# Extend the Y/N data validation list down to the new row
$ws.Range("D2:D28").Validation.Delete()
$ws.Range("D2:D28").Validation.Add(3, 1, 1, '"Y,N"')

# Update the view to reflect where the user ended up (scrolled/selected)
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C28").Select()
